$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after the existing one; it becomes the active tab
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet2"

# Header / label cell
$newSheet.Range("A1").Value = "flair text"

# Two helper ratio formulas
$newSheet.Range("B2").Formula = "=5013439/7518510"
$newSheet.Range("B3").Formula = "=2505071/7518510"

# Table header row
$newSheet.Range("B5").Value = "Row"
$newSheet.Range("C5").Value = "post_id_unique_count"
$newSheet.Range("D5").Value = "posts_with_flair_text"
$newSheet.Range("E5").Value = "flair_text_unique_count"

# Table data row
$newSheet.Range("B6").Value = 1
$newSheet.Range("C6").Value = 7518510
$newSheet.Range("D6").Value = 2505071
$newSheet.Range("E6").Value = 96734

# Match the thousands-separator number format already used elsewhere in the workbook
$newSheet.Range("C6:E6").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# Column widths to match the authored layout (stored width = ColumnWidth + 0.83)
$newSheet.Columns("C:D").ColumnWidth = 12.17
$newSheet.Columns("E:E").ColumnWidth = 10.17

# Selection/active cell matches what was captured when the sheet was saved
$excel.Goto($newSheet.Range("C5:E6"))
